$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Costes Humanos")
$ws2 = $wb.Worksheets.Item("Costes Materiales")
$ws3 = $wb.Worksheets.Item("Costes Proyecto")

# --- Sheet1: Costes Humanos ---
# Rename "Documentador" -> "Documentalist" everywhere it appears
$ws1.Range("D7").Value = "Documentalist"
$ws1.Range("F7").Value = "Documentalist"
$ws1.Range("B27").Value = "Documentalist"
$ws1.Range("B28").Value = "Documentalist"
$ws1.Range("B29").Value = "Documentalist"

# New salary data for Documentalist role (row 7 F:I block)
$ws1.Range("G7").Value = 31165
$ws1.Range("H7").Formula = "=G7/14"
$ws1.Range("I7").Formula = "=H7/30"

# Fix mislabeled role in D12 (was Solution Architect, should be Business Analyst)
$ws1.Range("D12").Value = "Business Analyst"

# Row 17 becomes the Documentalist cost row; row 18 becomes new TOTAL row
$ws1.Range("D17").Value = "Documentalist"
$ws1.Range("E17").Formula = "=PRODUCT(E7,I7)"
$ws1.Range("D18").Value = "TOTAL"
$ws1.Range("E18").Formula = "=SUM(E12:E17)"

Write-Output "sheet1 done"
